$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 15 & 16 (M3 screws / M3 nuts references) ---
$ws.Range("B15").Value = "M3 screws (12 mm)"
$ws.Range("D15").Value = "190-456"

$ws.Range("D16").Value = "189-563"

# --- Insert 3 new rows after row 16 for the new RS amidata parts ---
$ws.Rows("17:19").Insert()

$ws.Range("B17").Value = "M4 screws (20 mm)"
$ws.Range("C17").Value = "RS amidata"
$ws.Range("D17").Value = "227-6849"
$ws.Range("E17").Value = 1

$ws.Range("B18").Value = "M4 nuts"
$ws.Range("C18").Value = "RS amidata"
$ws.Range("D18").Value = "189-579"
$ws.Range("E18").Value = 1

$ws.Range("B19").Value = "Bosch Rexroth Connecting Component, Angle Bracket, strut profile 20 mm, groove Size 6mm"
$ws.Range("C19").Value = "RS amidata"
$ws.Range("D19").Value = "466-7354"
$ws.Range("E19").Value = 1
